$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "B-" + $cell.Text
}
